$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common / static values shared by all data rows in this block
$A = 1
$B = "Agrícola del Norte S.A. de Arica"
$C = "Arica y Parinacota"
$E = 15
$F = "Fruta"
$G = 100107
$H = "Otros"
$I = 100107011
$J = "Tuna"
$K = "Sin especificar"
$R = "Región de Coquimbo"

# Target data (rows 2-3 stay as-is; rows 4-9 are the new/changed block)
$rows = @(
    @{ Row = 4; D = 44679; L = "Segunda"; M = 200; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; S = 1475; T = 20 },
    @{ Row = 5; D = 44679; L = "Tercera"; M = 200; N = 24000; O = 25000; P = 24500; Q = "$/caja 20 kilos"; S = 1225; T = 20 },
    @{ Row = 6; D = 44664; L = "Segunda"; M = 150; N = 29000; O = 30000; P = 29500; Q = "$/caja 18 kilos"; S = 1639; T = 18 },
    @{ Row = 7; D = 44643; L = "Primera"; M = 160; N = 28000; O = 30000; P = 29000; Q = "$/caja 20 kilos"; S = 1450; T = 20 },
    @{ Row = 8; D = 44650; L = "Primera"; M = 160; N = 31000; O = 32000; P = 31500; Q = "$/caja 20 kilos"; S = 1575; T = 20 },
    @{ Row = 9; D = 44650; L = "Segunda"; M = 250; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; S = 1475; T = 20 }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Cells.Item($n, 1).Value = $A
    $ws.Cells.Item($n, 2).Value = $B
    $ws.Cells.Item($n, 3).Value = $C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($n, 5).Value = $E
    $ws.Cells.Item($n, 6).Value = $F
    $ws.Cells.Item($n, 7).Value = $G
    $ws.Cells.Item($n, 8).Value = $H
    $ws.Cells.Item($n, 9).Value = $I
    $ws.Cells.Item($n, 10).Value = $J
    $ws.Cells.Item($n, 11).Value = $K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
    $ws.Cells.Item($n, 14).Value = $r.N
    $ws.Cells.Item($n, 15).Value = $r.O
    $ws.Cells.Item($n, 16).Value = $r.P
    $ws.Cells.Item($n, 17).Value = $r.Q
    $ws.Cells.Item($n, 18).Value = $R
    $ws.Cells.Item($n, 19).Value = $r.S
    $ws.Cells.Item($n, 20).Value = $r.T
}
